# Add a new worksheet "L6" at the end of the workbook that summarises each
# team's last six games across Form / Goals scored / Goals conceded /
# Total Goals (alphabetical team order, matching the other summary sheets).

$wb = $excel.ActiveWorkbook

$headers = @(
    'Form',
    'Goals scored',
    'Goals conceded',
    'Total Goals'
)

$form = @(
    'Angers,D D L L L L',
    'Bordeaux,L L L L L W',
    'Brest,D L D D W L',
    'Dijon,L L L W L L',
    'Lens,W D W D W L',
    'Lille,L W W D W W',
    'Lorient,D W L L W W',
    'Lyon,L D W W L W',
    'Marseille,L W D W W D',
    'Metz,L L L D L W',
    'Monaco,W W W W W L',
    'Montpellier,W D D D L L',
    'Nantes,D L L L W W',
    'Nice,W W D L W L',
    'Nimes,W L D D L D',
    'Paris SG,W L W W W W',
    'Reims,W D D D L D',
    'Rennes,W D W W W L',
    'St Etienne,L W W L L W',
    'Strasbourg,L W L D L D'
)

$goalsScored = @(
    'Angers,0 1 0 0 0 0',
    'Bordeaux,1 2 1 0 1 1',
    'Brest,0 0 1 1 2 1',
    'Dijon,0 0 0 2 1 1',
    'Lens,2 1 4 1 2 1',
    'Lille,1 1 2 1 3 2',
    'Lorient,1 1 1 2 4 2',
    'Lyon,2 1 3 2 2 3',
    'Marseille,0 2 3 3 3 1',
    'Metz,1 0 0 0 1 5',
    'Monaco,4 4 3 3 1 2',
    'Montpellier,3 1 3 1 1 1',
    'Nantes,1 1 0 1 2 4',
    'Nice,3 2 0 0 3 0',
    'Nimes,2 0 1 1 1 2',
    'Paris SG,4 0 4 3 3 2',
    'Reims,1 2 0 0 1 2',
    'Rennes,3 2 1 3 5 0',
    'St Etienne,0 2 4 2 1 2',
    'Strasbourg,1 3 1 1 1 1'
)

$goalsConceded = @(
    'Angers,0 1 3 3 1 2',
    'Bordeaux,3 3 4 3 4 0',
    'Brest,0 1 1 1 1 4',
    'Dijon,1 2 3 0 5 5',
    'Lens,1 1 1 1 1 2',
    'Lille,2 0 0 1 2 0',
    'Lorient,1 0 4 3 1 0',
    'Lyon,4 1 0 1 3 2',
    'Marseille,3 0 3 2 1 1',
    'Metz,3 4 2 0 3 1',
    'Monaco,0 0 0 0 0 3',
    'Montpellier,1 1 3 1 3 2',
    'Nantes,1 2 1 2 1 1',
    'Nice,0 1 0 2 1 2',
    'Nimes,1 2 1 1 2 2',
    'Paris SG,2 1 1 2 1 1',
    'Reims,0 2 0 0 3 2',
    'Rennes,1 2 0 0 1 1',
    'St Etienne,4 0 1 3 2 1',
    'Strasbourg,2 2 4 1 2 1'
)

$totalGoals = @(
    'Angers,0 2 3 3 1 2',
    'Bordeaux,4 5 5 3 5 1',
    'Brest,0 1 2 2 3 5',
    'Dijon,1 2 3 2 6 6',
    'Lens,3 2 5 2 3 3',
    'Lille,3 1 2 2 5 2',
    'Lorient,2 1 5 5 5 2',
    'Lyon,6 2 3 3 5 5',
    'Marseille,3 2 6 5 4 2',
    'Metz,4 4 2 0 4 6',
    'Monaco,4 4 3 3 1 5',
    'Montpellier,4 2 6 2 4 3',
    'Nantes,2 3 1 3 3 5',
    'Nice,3 3 0 2 4 2',
    'Nimes,3 2 2 2 3 4',
    'Paris SG,6 1 5 5 4 3',
    'Reims,1 4 0 0 4 4',
    'Rennes,4 4 1 3 6 1',
    'St Etienne,4 2 5 5 3 3',
    'Strasbourg,3 5 5 2 3 2'
)

$tableWs = $wb.Worksheets.Item("Table")

# The Table sheet's A1 was only ever an empty placeholder cell; drop it so
# the saved sheet1.xml no longer carries a pointless blank <c r="A1"/>.
$tableWs.Range("A1").ClearContents() | Out-Null

# New sheet goes after the last existing sheet ("Goal totals v2").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "L6"

# Header row (B1:E1)
$ws.Range("B1").Value = $headers[0]
$ws.Range("C1").Value = $headers[1]
$ws.Range("D1").Value = $headers[2]
$ws.Range("E1").Value = $headers[3]

# Column A (rows 2-21) reuses the "1".."20" row-number labels already on
# the Table sheet, so copy them across instead of re-typing (keeps them as
# the same shared-string text cells rather than turning into numbers).
$tableWs.Range("A2:A21").Copy() | Out-Null
$ws.Range("A2").PasteSpecial() | Out-Null
$excel.CutCopyMode = $false

# Fill column-by-column (all of "Form" first, then "Goals scored", etc.)
# so new shared-string entries are interned in the same order as the
# original author's save.
for ($i = 0; $i -lt 20; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $form[$i]
}
for ($i = 0; $i -lt 20; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $goalsScored[$i]
}
for ($i = 0; $i -lt 20; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $goalsConceded[$i]
}
for ($i = 0; $i -lt 20; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $totalGoals[$i]
}

# Column widths for B/C/D to fit the longer "Team,results" strings.
$ws.Columns.Item(2).ColumnWidth = 22.0
$ws.Columns.Item(3).ColumnWidth = 22.6
$ws.Columns.Item(4).ColumnWidth = 26.0

# Make L6 the active sheet / selected cell, matching the saved view state.
$ws.Select() | Out-Null
$ws.Range("C14").Select() | Out-Null
